$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typos in "Contact form" (was "Concact form" / "Contact orm")
$ws.Range("D4").Value = "Contact form"
$ws.Range("D8").Value = "Contact form"
$ws.Range("D9").Value = "Contact form"

# Row 32 (Graftek engineering / Indeed): mark response received, interview = no
$ws.Range("F32").Value = "x"
$ws.Range("G32").Value = "n"

# Row 36 (ARCTEC Solutions / Indeed): mark response received, interview = yes, add notes
$ws.Range("F36").Value = "x"
$ws.Range("G36").Value = "y"
$ws.Range("I36").Value = "Radio wave, ham radio"

# Update the active selection to match the saved view state
$ws.Range("B19").Select()
